$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibition) and "全部类型" (All Types) sheets list the
# same events in rows 2-6 and need the same F-column (想去人数 / interest
# count) updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5760
    $ws.Range("F3").Value = 11
    $ws.Range("F5").Value = 969
    $ws.Range("F6").Value = 61
}
